$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Formula = "'29.349.43"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.Formula = "'  +0.44%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.Formula = "'1.843.14"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.Formula = "'  +0.12%  "
$c.ClearFormats()

$c = $ws.Range("D4")
$c.Formula = "'0.9985"
$c.ClearFormats()
$c = $ws.Range("E4")
$c.Formula = "'  -0.33%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.Formula = "'240.02"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.Formula = "'  -0.23%  "
$c.ClearFormats()

$c = $ws.Range("D6")
$c.Formula = "'0.6304"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.Formula = "'  +0.51%  "
$c.ClearFormats()

$c = $ws.Range("E7")
$c.Formula = "'  -0.33%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.Formula = "'0.07478"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.Formula = "'  +0.53%  "
$c.ClearFormats()

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range("D9")
$c.Formula = "'0.2900"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.Formula = "'  +0.52%  "
$c.ClearFormats()

$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range("D10")
$c.Formula = "'24.97"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.Formula = "'  +3.11%  "
$c.ClearFormats()

$c = $ws.Range("D11")
$c.Formula = "'0.07730"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.Formula = "'  -0.02%  "
$c.ClearFormats()

$c = $ws.Range("D12")
$c.Formula = "'1.847.15"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.Formula = "'  +0.34%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.Formula = "'4.979"
$c.ClearFormats()
$c = $ws.Range("E13")
$c.Formula = "'  -0.03%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.Formula = "'0.6774"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.Formula = "'  +0.31%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.Formula = "'0.00001035"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.Formula = "'  +2.71%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.Formula = "'81.92"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.Formula = "'  -0.11%  "
$c.ClearFormats()

$c = $ws.Range("D17")
$c.Formula = "'6.230"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.Formula = "'  +1.98%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.Formula = "'29.373.46"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.Formula = "'  +0.34%  "
$c.ClearFormats()

$c = $ws.Range("D19")
$c.Formula = "'229.16"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.Formula = "'  +0.91%  "
$c.ClearFormats()

$c = $ws.Range("E20")
$c.Formula = "'  +0.58%  "
$c.ClearFormats()

$c = $ws.Range("D21")
$c.Formula = "'0.9999"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.Formula = "'  -0.29%  "
$c.ClearFormats()

$c = $ws.Range("D22")
$c.Formula = "'7.389"
$c.ClearFormats()
$c = $ws.Range("E22")
$c.Formula = "'  +0.65%  "
$c.ClearFormats()

$c = $ws.Range("D23")
$c.Formula = "'0.9996"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.Formula = "'  -0.34%  "
$c.ClearFormats()

$c = $ws.Range("D24")
$c.Formula = "'158.13"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.Formula = "'  -0.39%  "
$c.ClearFormats()

$c = $ws.Range("D25")
$c.Formula = "'8.530"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.Formula = "'  +1.91%  "
$c.ClearFormats()

$c = $ws.Range("D26")
$c.Formula = "'0.1356"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.Formula = "'  -1.00%  "
$c.ClearFormats()

$c = $ws.Range("D27")
$c.Formula = "'17.48"
$c.ClearFormats()
$c = $ws.Range("E27")
$c.Formula = "'  -0.23%  "
$c.ClearFormats()

$c = $ws.Range("D28")
$c.Formula = "'0.06859"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.Formula = "'  +12.20%  "
$c.ClearFormats()

$c = $ws.Range("D29")
$c.Formula = "'1.455"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.Formula = "'  +4.29%  "
$c.ClearFormats()

$c = $ws.Range("D30")
$c.Formula = "'1.488"
$c.ClearFormats()

$c = $ws.Range("D31")
$c.Formula = "'4.064"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.Formula = "'  -0.13%  "
$c.ClearFormats()

$c = $ws.Range("D32")
$c.Formula = "'4.066"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.Formula = "'  +0.89%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.Formula = "'1.833"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.Formula = "'  +0.91%  "
$c.ClearFormats()

$c = $ws.Range("D35")
$c.Formula = "'0.7005"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.Formula = "'  +0.72%  "
$c.ClearFormats()

$c = $ws.Range("D36")
$c.Formula = "'2.585"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.Formula = "'  -0.35%  "
$c.ClearFormats()

$c = $ws.Range("D37")
$c.Formula = "'0.01846"
$c.ClearFormats()
$c = $ws.Range("E37")
$c.Formula = "'  +1.99%  "
$c.ClearFormats()

$c = $ws.Range("D38")
$c.Formula = "'2.820"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.Formula = "'  -0.12%  "
$c.ClearFormats()

$c = $ws.Range("D39")
$c.Formula = "'1.238.73"
$c.ClearFormats()
$c = $ws.Range("E39")
$c.Formula = "'  -0.24%  "
$c.ClearFormats()

$c = $ws.Range("D40")
$c.Formula = "'6.773"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.Formula = "'  +4.42%  "
$c.ClearFormats()

$c = $ws.Range("D41")
$c.Formula = "'0.9443"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.Formula = "'  +4.19%  "
$c.ClearFormats()

$c = $ws.Range("D42")
$c.Formula = "'0.9990"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.Formula = "'  -0.10%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.Formula = "'2.002.83"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.Formula = "'  +0.16%  "
$c.ClearFormats()

$c = $ws.Range("D44")
$c.Formula = "'101.01"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.Formula = "'  -0.46%  "
$c.ClearFormats()

$c = $ws.Range("D45")
$c.Formula = "'65.56"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.Formula = "'  -0.66%  "
$c.ClearFormats()

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D46")
$c.Formula = "'7.051"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.Formula = "'  +0.37%  "
$c.ClearFormats()

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D47")
$c.Formula = "'1.717"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.Formula = "'  +3.99%  "
$c.ClearFormats()

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range("D48")
$c.Formula = "'8.976"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.Formula = "'  +0.03%  "
$c.ClearFormats()

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D49")
$c.Formula = "'0.1146"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.Formula = "'  -1.20%  "
$c.ClearFormats()

$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D50")
$c.Formula = "'0.3917"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.Formula = "'  -0.36%  "
$c.ClearFormats()

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D51")
$c.Formula = "'0.05677"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.Formula = "'  -0.30%  "
$c.ClearFormats()
